# (#33) Alteração nos rótulos da tabela para já transformar a primeira linha
# em cabeçalho automaticamente no power bi
#
# For each worksheet, prefix the year/interval labels on row 1 (B1, C1, ...)
# with "Ano " (for plain years) or "Intervalo " (for year ranges), so the
# labels become self-describing headers (e.g. "2015" -> "Ano 2015",
# "2015-2030" -> "Intervalo 2015-2030").

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $lastCol = $used.Columns.Count

    for ($col = 2; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item(1, $col)
        $val = $cell.Value2

        if ($null -eq $val) { continue }

        $text = [string]$val

        if ($text -match "^\d{4}-\d{4}$") {
            $cell.Value2 = "Intervalo " + $text
        }
        elseif ($text -match "^\d{4}$") {
            $cell.Value2 = "Ano " + $text
        }
    }
}
